$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.969.80"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "2.350.63"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.60"
$ws.Range("E5").Value = "  -3.66%  "
$ws.Range("D6").Value = "143.75"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "2.353.62"
$ws.Range("E9").Value = "  -5.83%  "
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("E11").Value = "  -6.63%  "
$ws.Range("D12").Value = "0.318"
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "2.760.94"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "55.022.66"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("E16").Value = "  -5.69%  "
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("D18").Value = "2.351.69"
$ws.Range("E18").Value = "  -5.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "311.72"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("E21").Value = "  -5.52%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.61"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "55.87"
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("D28").Value = "2.453.35"
$ws.Range("E28").Value = "  -5.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("E29").Value = "  -6.18%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "0.0₃0751"
$ws.Range("E31").Value = "  -4.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.68"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.96"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("E36").Value = "  -5.42%  "
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.820"
$ws.Range("E38").Value = "  -4.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.45"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.34"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0944"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.576"
$ws.Range("E44").Value = "  -5.78%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0524"
$ws.Range("E45").Value = "  -6.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.14"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "254.13"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("E49").Value = "  -7.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.76"
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("D51").Value = "1.777.02"
$ws.Range("E51").Value = "  -6.23%  "
